$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing row 180, shifting rows 180-183 down to 181-184.
$ws.Rows.Item(180).Insert()

# Populate the newly inserted row 180 with the new weekly data point.
$ws.Cells.Item(180, 1).Value = 5
$ws.Cells.Item(180, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(180, 3).Value = "Maule"
$ws.Cells.Item(180, 4).Value = 44595
$ws.Cells.Item(180, 5).Value = 7
$ws.Cells.Item(180, 6).Value = 100112024
$ws.Cells.Item(180, 7).Value = "Choclo"
$ws.Cells.Item(180, 8).Value = "Choclero"
$ws.Cells.Item(180, 9).Value = "Primera"
$ws.Cells.Item(180, 10).Value = 50000
$ws.Cells.Item(180, 11).Value = 100
$ws.Cells.Item(180, 12).Value = 100
$ws.Cells.Item(180, 13).Value = 100
$ws.Cells.Item(180, 14).Value = '$/unidad'
$ws.Cells.Item(180, 15).Value = "Región del Maule"
$ws.Cells.Item(180, 16).Value = 100
$ws.Cells.Item(180, 17).Value = 1
$ws.Cells.Item(180, 18).Value = "Hortaliza"

# Match the date number format already used by the other date cells in column D.
$ws.Cells.Item(180, 4).NumberFormat = $ws.Cells.Item(181, 4).NumberFormat
